# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new value for column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13821
$ws1.Range("F7").Value = 250
$ws1.Range("F8").Value = 1777
$ws1.Range("F10").Value = 134
$ws1.Range("F13").Value = 531
$ws1.Range("F16").Value = 13869
$ws1.Range("F19").Value = 14940
$ws1.Range("F21").Value = 8228
$ws1.Range("F30").Value = 1034
$ws1.Range("F37").Value = 216

# Sheet "全部类型": row -> new value for column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13821
$ws4.Range("F7").Value = 250
$ws4.Range("F8").Value = 1777
$ws4.Range("F10").Value = 134
$ws4.Range("F13").Value = 531
$ws4.Range("F16").Value = 13869
$ws4.Range("F19").Value = 14940
$ws4.Range("F21").Value = 8228
$ws4.Range("F30").Value = 1034
$ws4.Range("F39").Value = 216
